# Apply cryptos list price/volume updates (commit: "Updated cryptos list on Wed Nov 15 16:15:53 UTC 2023 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '36.387.01'
$ws.Cells.Item(2, 5).Value = '  +0.15%  '
$ws.Cells.Item(3, 4).Value = '2.017.17'
$ws.Cells.Item(3, 5).Value = '  -1.50%  '
$ws.Cells.Item(4, 5).Value = '  +0.08%  '
$ws.Cells.Item(5, 4).Value = '''251.41'
$ws.Cells.Item(5, 5).Value = '  +2.63%  '
$ws.Cells.Item(6, 5).Value = '  -2.83%  '
$ws.Cells.Item(7, 4).Value = '''62.03'
$ws.Cells.Item(7, 5).Value = '  +9.29%  '
$ws.Cells.Item(8, 5).Value = '  +0.06%  '
$ws.Cells.Item(9, 4).Value = '''59.17'
$ws.Cells.Item(9, 5).Value = '  -8.59%  '
$ws.Cells.Item(10, 5).Value = '  +0.70%  '
$ws.Cells.Item(11, 5).Value = '  -0.33%  '
$ws.Cells.Item(12, 5).Value = '  -1.63%  '
$ws.Cells.Item(13, 4).Value = '''0.907'
$ws.Cells.Item(13, 5).Value = '  -0.89%  '
$ws.Cells.Item(14, 4).Value = '''14.82'
$ws.Cells.Item(14, 5).Value = '  +3.95%  '
$ws.Cells.Item(15, 4).Value = '2.313.87'
$ws.Cells.Item(16, 4).Value = '''20.34'
$ws.Cells.Item(16, 5).Value = '  +14.74%  '
$ws.Cells.Item(17, 5).Value = '  +0.62%  '
$ws.Cells.Item(18, 4).Value = '2.022.16'
$ws.Cells.Item(18, 5).Value = '  -0.86%  '
$ws.Cells.Item(19, 4).Value = '36.387.32'
$ws.Cells.Item(19, 5).Value = '  +0.36%  '
$ws.Cells.Item(20, 4).Value = '''72.03'
$ws.Cells.Item(20, 5).Value = '  +0.85%  '
$ws.Cells.Item(21, 5).Value = '  +0.79%  '
$ws.Cells.Item(22, 4).Value = '''5.30'
$ws.Cells.Item(22, 5).Value = '  +1.95%  '
$ws.Cells.Item(23, 4).Value = '''234.34'
$ws.Cells.Item(23, 5).Value = '  -1.05%  '
$ws.Cells.Item(24, 5).Value = '  +17.59%  '
$ws.Cells.Item(25, 5).Value = '  -0.17%  '
$ws.Cells.Item(26, 4).Value = '''2.32'
$ws.Cells.Item(26, 5).Value = '  -1.42%  '
$ws.Cells.Item(27, 4).Value = '''9.59'
$ws.Cells.Item(27, 5).Value = '  +3.20%  '
$ws.Cells.Item(28, 4).Value = '''163.59'
$ws.Cells.Item(28, 5).Value = '  -0.71%  '
$ws.Cells.Item(29, 4).Value = '''19.61'
$ws.Cells.Item(29, 5).Value = '  -1.81%  '
$ws.Cells.Item(30, 2).Value = 'Stellar'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(30, 4).Value = '''0.120'
$ws.Cells.Item(30, 5).Value = '  -0.70%  '
$ws.Cells.Item(31, 2).Value = 'Filecoin'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(31, 4).Value = '''5.10'
$ws.Cells.Item(31, 5).Value = '  +2.46%  '
$ws.Cells.Item(32, 2).Value = 'Kaspa'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(32, 4).Value = '''0.111'
$ws.Cells.Item(32, 5).Value = '  +27.58%  '
$ws.Cells.Item(33, 5).Value = '  -1.51%  '
$ws.Cells.Item(34, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(34, 4).Value = '''4.59'
$ws.Cells.Item(34, 5).Value = '  +3.85%  '
$ws.Cells.Item(35, 2).Value = 'Hedera'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(35, 4).Value = '''0.0608'
$ws.Cells.Item(35, 5).Value = '  +1.53%  '
$ws.Cells.Item(36, 5).Value = '  +11.03%  '
$ws.Cells.Item(37, 5).Value = '  +0.12%  '
$ws.Cells.Item(38, 5).Value = '  -0.59%  '
$ws.Cells.Item(39, 4).Value = '''5.88'
$ws.Cells.Item(39, 5).Value = '  +16.75%  '
$ws.Cells.Item(40, 5).Value = '  +14.87%  '
$ws.Cells.Item(41, 5).Value = '  +0.66%  '
$ws.Cells.Item(42, 5).Value = '  +1.98%  '
$ws.Cells.Item(43, 2).Value = 'ARBITRUM'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(43, 4).Value = '''1.13'
$ws.Cells.Item(43, 5).Value = '  +2.71%  '
$ws.Cells.Item(44, 2).Value = 'VeChain'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(44, 4).Value = '''0.0216'
$ws.Cells.Item(44, 5).Value = '  +0.40%  '
$ws.Cells.Item(45, 2).Value = 'Maker'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(45, 4).Value = '1.447.95'
$ws.Cells.Item(45, 5).Value = '  +5.60%  '
$ws.Cells.Item(46, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(46, 4).Value = '''16.68'
$ws.Cells.Item(46, 5).Value = '  +5.02%  '
$ws.Cells.Item(47, 2).Value = 'Aave'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(47, 4).Value = '''94.66'
$ws.Cells.Item(47, 5).Value = '  +1.10%  '
$ws.Cells.Item(48, 2).Value = 'FraxShare'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(48, 4).Value = '''7.87'
$ws.Cells.Item(48, 5).Value = '  +6.20%  '
$ws.Cells.Item(49, 4).Value = '''2.62'
$ws.Cells.Item(49, 5).Value = '  +14.42%  '
$ws.Cells.Item(50, 4).Value = '''2.95'
$ws.Cells.Item(50, 5).Value = '  +0.24%  '
$ws.Cells.Item(51, 4).Value = '''47.27'
$ws.Cells.Item(51, 5).Value = '  +3.12%  '
